$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column (C) for rows 2-28
# from serial date 45454 (2024-06-11) to 45455 (2024-06-12)
for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45454) {
        $cell.Value2 = 45455
    }
}
